# ============================================================
# Add 2022-Q4 data
# ============================================================
$wb = $excel.ActiveWorkbook

# ---- Sheet shuffle -----------------------------------------
# The existing "2022-Q3" sheet (sheetId=2) is duplicated. The
# duplicate keeps the old Q3 detail data untouched and becomes the
# new "2022-Q3" sheet (sheetId=3). The original is renamed to
# "2022-Q4" (keeps sheetId=2) and repopulated with new data below.
$q3src = $wb.Worksheets.Item("2022-Q3")
$q3src.Copy($null, $q3src)
$q3dup = $wb.Worksheets.Item($q3src.Index + 1)
$q3src.Name = "2022-Q4"
$q3dup.Name = "2022-Q3"
$q4 = $q3src

# ---- Clear old content out of the renamed "2022-Q4" sheet ---
$q4.Cells.Clear()

# ---- Borrow header / index-column formatting from "总计" -----
$summary = $wb.Worksheets.Item("总计")
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q4.Range("A2:A13").PasteSpecial(-4122)

# ---- Header row ----------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# ---- Format the text columns (B:G) as Text so codes / decimal
# strings like "010709" or "27.17" are NOT auto-coerced to numbers
$q4.Range("B2:G13").NumberFormat = "@"

# ---- Data rows -------------------------------------------------
# row 2
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "010709"
$q4.Range("C2").Value = "安信医药健康主题股票A"
$q4.Range("D2").Value = "27.17"
$q4.Range("E2").Value = "94.69"
$q4.Range("F2").Value = "4.38"
$q4.Range("G2").Value = "1.1900"
$q4.Range("H2").Value = 7
# row 3
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "010710"
$q4.Range("C3").Value = "安信医药健康主题股票C"
$q4.Range("D3").Value = "24.37"
$q4.Range("E3").Value = "94.69"
$q4.Range("F3").Value = "4.38"
$q4.Range("G3").Value = "1.0674"
$q4.Range("H3").Value = 7
# row 4
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "160926"
$q4.Range("C4").Value = "大成创业板两年定期开放混合A"
$q4.Range("D4").Value = "7.36"
$q4.Range("E4").Value = "80.56"
$q4.Range("F4").Value = "7.42"
$q4.Range("G4").Value = "0.5461"
$q4.Range("H4").Value = 2
# row 5
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "501079"
$q4.Range("C5").Value = "大成科创主题混合（LOF）A"
$q4.Range("D5").Value = "9.55"
$q4.Range("E5").Value = "80.68"
$q4.Range("F5").Value = "5.48"
$q4.Range("G5").Value = "0.5233"
$q4.Range("H5").Value = 8
# row 6
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "012473"
$q4.Range("C6").Value = "大成成长回报六个月持有混合A"
$q4.Range("D6").Value = "7.49"
$q4.Range("E6").Value = "75.28"
$q4.Range("F6").Value = "5.33"
$q4.Range("G6").Value = "0.3992"
$q4.Range("H6").Value = 5
# row 7
$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "010371"
$q4.Range("C7").Value = "大成成长进取混合A"
$q4.Range("D7").Value = "3.61"
$q4.Range("E7").Value = "74.21"
$q4.Range("F7").Value = "5.36"
$q4.Range("G7").Value = "0.1935"
$q4.Range("H7").Value = 3
# row 8
$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "009798"
$q4.Range("C8").Value = "大成创业板两年定期开放混合C"
$q4.Range("D8").Value = "2.60"
$q4.Range("E8").Value = "80.56"
$q4.Range("F8").Value = "7.42"
$q4.Range("G8").Value = "0.1929"
$q4.Range("H8").Value = 2
# row 9
$q4.Range("A9").Value = 7
$q4.Range("B9").Value = "010372"
$q4.Range("C9").Value = "大成成长进取混合C"
$q4.Range("D9").Value = "1.50"
$q4.Range("E9").Value = "74.21"
$q4.Range("F9").Value = "5.36"
$q4.Range("G9").Value = "0.0804"
$q4.Range("H9").Value = 3
# row 10
$q4.Range("A10").Value = 8
$q4.Range("B10").Value = "005652"
$q4.Range("C10").Value = "国富天颐混合A"
$q4.Range("D10").Value = "5.41"
$q4.Range("E10").Value = "21.77"
$q4.Range("F10").Value = "0.60"
$q4.Range("G10").Value = "0.0325"
$q4.Range("H10").Value = 8
# row 11
$q4.Range("A11").Value = 9
$q4.Range("B11").Value = "012474"
$q4.Range("C11").Value = "大成成长回报六个月持有混合C"
$q4.Range("D11").Value = "0.37"
$q4.Range("E11").Value = "75.28"
$q4.Range("F11").Value = "5.33"
$q4.Range("G11").Value = "0.0197"
$q4.Range("H11").Value = 5
# row 12
$q4.Range("A12").Value = 10
$q4.Range("B12").Value = "005653"
$q4.Range("C12").Value = "国富天颐混合C"
$q4.Range("D12").Value = "0.60"
$q4.Range("E12").Value = "21.77"
$q4.Range("F12").Value = "0.60"
$q4.Range("G12").Value = "0.0036"
$q4.Range("H12").Value = 8
# row 13
$q4.Range("A13").Value = 11
$q4.Range("B13").Value = "016198"
$q4.Range("C13").Value = "大成科创主题混合（LOF）C"
$q4.Range("D13").Value = "0.06"
$q4.Range("E13").Value = "80.68"
$q4.Range("F13").Value = "5.48"
$q4.Range("G13").Value = "0.0033"
$q4.Range("H13").Value = 8

# ---- Update the "总计" summary sheet --------------------------
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 12
$summary.Range("D2").Value = 4.25

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 13
$summary.Range("D3").Value = 2.34

Write-Host "done"
